$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy header style (from B1, which already has the bold/border/center style) to the new header cells J1:W1
$ws.Range("B1").Copy() | Out-Null
$ws.Range("J1:W1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Step 2: write header row (A1 unchanged = "Date")
$headers = @("Date", "Alexis Rainey", "Balduzzi", "Burns", "Curley", "Doyle", "Espona", "Ferriolo", "Hackman", "Holzman", "Hughes", "Johnson", "McCann", "McFadden", "Medico", "Myers", "Pla", "Reilly", "Rodrigo", "Streib", "Tollaksen", "Wasyliw", "Yanovich")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Step 3: write the data rows 2-6 (column A / dates are unchanged)
$row2 = @(665.2515, 640.1744, 670.8605, "", 779.6707, 824.9546, 603.0765, 695.3825, 658.4766, "", 725.4476, 790.4619, 788.4167, 677.2832, 715.1943, 663.6866, "", 661.7878, "", 800.2089, 272.5958, 637.4342)
for ($i = 0; $i -lt $row2.Count; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}

$row3 = @(519.0682, 420.1565, 298.4097, 294.0178, 461.3162, 452.4021, 279.6842, 408.1442, 433.3555, 288.2412, 603.7235, 545.6038, 440.5764, 535.363, 360.5946, 356.5623, 454.0465, 498.1062, "", 445.3439, 146.4231, 420.0008)
for ($i = 0; $i -lt $row3.Count; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
}

$row4 = @(540.2351, 456.9402, 478.8517, 260.6904, 665.1731, 518.8987, 624.0218, 634.8641, 602.8357, 445.3059, 703.6967, 703.6943, 578.3936, 770.2661, 579.48, 520.8793, 569.0826, 567.2286, "", 696.1316, 224.0441, "")
for ($i = 0; $i -lt $row4.Count; $i++) {
    $ws.Cells.Item(4, $i + 2).Value = $row4[$i]
}

$row5 = @(497.8744, 457.7425, "", 269.8513, 578.1309, 596.7133, 383.2716, 435.0821, 425.721, 373.8116, 325.4119, 528.1039, 545.8407, 508.8996, "", 481.3024, 438.782, 491.3095, 211.6593, 554.0063, 191.6313, "")
for ($i = 0; $i -lt $row5.Count; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $row5[$i]
}

$row6 = @(502.4344, 393.6091, "", 187.3585, 478.5848, 515.1994, 393.3857, 489.3714, 275.5775, 157.8899, 484.0656, 261.3154, 348.7925, 517.1587, "", 356.5423, 390.0852, 356.0464, 182.1284, "", 136.9055, "")
for ($i = 0; $i -lt $row6.Count; $i++) {
    $ws.Cells.Item(6, $i + 2).Value = $row6[$i]
}
